$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D and E columns for the affected rows so that
# numeric-looking strings (e.g. "582.19", "0.519") are stored as text,
# matching the original inlineStr cell type instead of being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '62.353.71'
$ws.Range("E2").Value = '  -2.14%  '

$ws.Range("D3").Value = '3.009.11'
$ws.Range("E3").Value = '  -2.19%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '582.19'
$ws.Range("E5").Value = '  -0.88%  '

$ws.Range("D6").Value = '146.92'
$ws.Range("E6").Value = '  -5.08%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = '3.010.55'
$ws.Range("E8").Value = '  -2.10%  '

$ws.Range("D9").Value = '0.519'
$ws.Range("E9").Value = '  -3.36%  '

$ws.Range("D10").Value = '0.148'
$ws.Range("E10").Value = '  -5.18%  '

$ws.Range("D11").Value = '5.64'
$ws.Range("E11").Value = '  -3.54%  '

$ws.Range("D12").Value = '0.440'
$ws.Range("E12").Value = '  -2.13%  '

$ws.Range("D13").Value = '0.0000228'
$ws.Range("E13").Value = '  -3.90%  '

$ws.Range("D14").Value = '34.70'
$ws.Range("E14").Value = '  -5.51%  '

$ws.Range("D15").Value = '0.121'
$ws.Range("E15").Value = '  +1.84%  '

$ws.Range("D16").Value = '3.508.29'
$ws.Range("E16").Value = '  -2.06%  '

$ws.Range("D17").Value = '7.03'
$ws.Range("E17").Value = '  -1.82%  '

$ws.Range("D18").Value = '62.342.39'
$ws.Range("E18").Value = '  -1.97%  '

$ws.Range("D19").Value = '3.011.06'
$ws.Range("E19").Value = '  -2.09%  '

$ws.Range("D20").Value = '456.90'
$ws.Range("E20").Value = '  -3.10%  '

$ws.Range("D21").Value = '13.89'
$ws.Range("E21").Value = '  -2.85%  '

$ws.Range("D22").Value = '0.680'
$ws.Range("E22").Value = '  -3.32%  '

$ws.Range("D23").Value = '7.31'
$ws.Range("E23").Value = '  -2.78%  '

$ws.Range("E24").Value = '  -6.19%  '

$ws.Range("D25").Value = '79.96'
$ws.Range("E25").Value = '  -0.54%  '

$ws.Range("D26").Value = '12.30'
$ws.Range("E26").Value = '  -4.17%  '

$ws.Range("D27").Value = '10.11'
$ws.Range("E27").Value = '  -2.73%  '

$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.34%  '

$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").Value = '7.17'
$ws.Range("E30").Value = '  -2.58%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '2.61'
$ws.Range("E31").Value = '  -1.62%  '

$ws.Range("D32").Value = '2.09'
$ws.Range("E32").Value = '  -1.97%  '

$ws.Range("D33").Value = '26.93'
$ws.Range("E33").Value = '  -0.74%  '

$ws.Range("D34").Value = '0.107'
$ws.Range("E34").Value = '  -5.14%  '

$ws.Range("E35").Value = '  -1.69%  '

$ws.Range("D36").Value = '0.0₃0790'
$ws.Range("E36").Value = '  -4.53%  '

$ws.Range("D37").Value = '5.73'
$ws.Range("E37").Value = '  -4.26%  '

$ws.Range("D38").Value = '2.12'
$ws.Range("E38").Value = '  -3.97%  '

$ws.Range("D39").Value = '50.35'
$ws.Range("E39").Value = '  -0.48%  '

$ws.Range("D40").Value = '8.98'
$ws.Range("E40").Value = '  -1.80%  '

$ws.Range("D41").Value = '2.88'
$ws.Range("E41").Value = '  -11.24%  '

$ws.Range("D42").Value = '416.38'
$ws.Range("E42").Value = '  -4.47%  '

$ws.Range("E43").Value = '  +0.79%  '

$ws.Range("E44").Value = '  -5.19%  '

$ws.Range("E45").Value = '  -1.81%  '

$ws.Range("D46").Value = '2.770.09'
$ws.Range("E46").Value = '  -1.34%  '

$ws.Range("D47").Value = '37.98'
$ws.Range("E47").Value = '  -4.83%  '

$ws.Range("D48").Value = '128.67'
$ws.Range("E48").Value = '  -1.14%  '

$ws.Range("E50").Value = '  -1.26%  '

$ws.Range("D51").Value = '23.76'
$ws.Range("E51").Value = '  -5.00%  '

# Restore original (default) style so no spurious style index is left on cells
$ws.Range("D2:E51").Style = "Normal"